{"js": "// Fix four Korean marketing-campaign bullet paragraphs so the\n// \"Munson's: ...\" taglines/slogans/mottos/phrases are rendered in English\n// (matching the quoted English tagline used elsewhere) and the Korean\n// sentences around them read naturally, with \"Munson's\" (not \"Munson\uc758\"/\n// \"\ubb38\uc2a8\uc758\") used consistently as the possessive brand reference.\n\nconst replacements = [\n  {\n    old: \"\ub9c8\ucf00\ud305 \ucea0\ud398\uc778\uc740 \\\"Munson's: Pickles and Preserves with a Purpose\\\"\ub77c\ub294 \ud0dc\uadf8\ub77c\uc778\uc744 \uc0ac\uc6a9\ud558\uc5ec Munson\uc758 \ube0c\ub79c\ub4dc \ubcf8\uc9c8\uc744 \ud3ec\ucc29\ud569\ub2c8\ub2e4.\",\n    new: \"\ub9c8\ucf00\ud305 \ucea0\ud398\uc778\uc740 \\\"Munson's: Pickles and Preserves with a Purpose\\\"\ub77c\ub294 \ud0dc\uadf8\ub77c\uc778\uc744 \uc0ac\uc6a9\ud558\uc5ec Munson's\uc758 \ube0c\ub79c\ub4dc \ubcf8\uc9c8\uc744 \ud3ec\ucc29\ud569\ub2c8\ub2e4.\"\n  },\n  {\n    old: \"\ub9c8\ucf00\ud305 \ucea0\ud398\uc778\uc740 \\\"Munson's: Just Pickles and Preserves \uc774\uc0c1\\\"\uc774\ub77c\ub294 \ubb38\uc2a8\uc758 \uc81c\ud488 \ud61c\ud0dd\uc744 \uac15\uc870\ud558\uae30 \uc704\ud574 \ub2e4\uc74c \uc2ac\ub85c\uac74\uc744 \uc0ac\uc6a9\ud569\ub2c8\ub2e4.\",\n    new: \"\ub9c8\ucf00\ud305 \ucea0\ud398\uc778\uc740 \\\"Munson's: More than Just Pickles and Preserves\\\"\ub77c\ub294 \uc2ac\ub85c\uac74\uc744 \uc0ac\uc6a9\ud558\uc5ec Munson's \uc81c\ud488\uc758 \uc774\uc810\uc744 \uac15\uc870\ud569\ub2c8\ub2e4.\"\n  },\n  {\n    old: \"\ub9c8\ucf00\ud305 \ucea0\ud398\uc778\uc740 \ub2e4\uc74c \ubaa8\ud1a0\ub97c \uc0ac\uc6a9\ud558\uc5ec Munson\uc758 \uace0\uac1d \uc639\ud638\uc5d0 \uc601\uac10\uc744 \uc90d\ub2c8\ub2e4. \\\"Munson's: \ud53c\ud074\uacfc \ubcf4\uc874\uc758 \uc0ac\ub791\uc744 \uacf5\uc720\ud558\uc2ed\uc2dc\uc624\\\".\",\n    new: \"\ub9c8\ucf00\ud305 \ucea0\ud398\uc778\uc740 \\\"Munson's: Share the Love of Pickles and Preserves\\\"\ub77c\ub294 \ubaa8\ud1a0\ub97c \uc0ac\uc6a9\ud558\uc5ec Munson's \uace0\uac1d\uc758 \uc9c0\uc9c0\ub3c4\uc5d0 \uc601\uac10\uc744 \uc90d\ub2c8\ub2e4.\"\n  },\n  {\n    old: \"\ub9c8\ucf00\ud305 \ucea0\ud398\uc778\uc740 \ub2e4\uc74c \ubb38\uad6c\ub97c \uc0ac\uc6a9\ud558\uc5ec Munson\uc758 \uc81c\ud488 \ud3c9\uac00\ud310\uc744 \uad6c\ub3d9\ud558\uace0 \uad6c\ub9e4\ud569\ub2c8\ub2e4. \\\"Munson's: \ucc3e\uae30, \uc0ac\uc6a9\ud574 \ubcf4\uae30, \uc0ac\ub791\ud558\uc138\uc694\\\".\",\n    new: \"\ub9c8\ucf00\ud305 \ucea0\ud398\uc778\uc740 \\\"Munson's: Find Them, Try Them, Love Them\\\"\ub77c\ub294 \ubb38\uad6c\ub97c \uc0ac\uc6a9\ud558\uc5ec Munson's \uc81c\ud488\uc744 \uccb4\ud5d8\ud558\uace0 \uad6c\ub9e4\ud558\ub3c4\ub85d \ud569\ub2c8\ub2e4.\"\n  }\n];\n\nconst body = context.document.body;\n\nfor (const { old, new: replacement } of replacements) {\n  const results = body.search(old, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + old);\n  }\n\n  for (const range of results.items) {\n    range.insertText(replacement, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Fix four Korean marketing-campaign bullet paragraphs so the\n# \"Munson's: ...\" taglines/slogans/mottos/phrases are rendered in English\n# (matching the quoted English tagline used elsewhere) and the Korean\n# sentences around them read naturally, with \"Munson's\" (not \"Munson\uc758\"/\n# \"\ubb38\uc2a8\uc758\") used consistently as the possessive brand reference.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @{\n    Old = \"\ub9c8\ucf00\ud305 \ucea0\ud398\uc778\uc740 `\"Munson's: Pickles and Preserves with a Purpose`\"\ub77c\ub294 \ud0dc\uadf8\ub77c\uc778\uc744 \uc0ac\uc6a9\ud558\uc5ec Munson\uc758 \ube0c\ub79c\ub4dc \ubcf8\uc9c8\uc744 \ud3ec\ucc29\ud569\ub2c8\ub2e4.\"\n    New = \"\ub9c8\ucf00\ud305 \ucea0\ud398\uc778\uc740 `\"Munson's: Pickles and Preserves with a Purpose`\"\ub77c\ub294 \ud0dc\uadf8\ub77c\uc778\uc744 \uc0ac\uc6a9\ud558\uc5ec Munson's\uc758 \ube0c\ub79c\ub4dc \ubcf8\uc9c8\uc744 \ud3ec\ucc29\ud569\ub2c8\ub2e4.\"\n  },\n  @{\n    Old = \"\ub9c8\ucf00\ud305 \ucea0\ud398\uc778\uc740 `\"Munson's: Just Pickles and Preserves \uc774\uc0c1`\"\uc774\ub77c\ub294 \ubb38\uc2a8\uc758 \uc81c\ud488 \ud61c\ud0dd\uc744 \uac15\uc870\ud558\uae30 \uc704\ud574 \ub2e4\uc74c \uc2ac\ub85c\uac74\uc744 \uc0ac\uc6a9\ud569\ub2c8\ub2e4.\"\n    New = \"\ub9c8\ucf00\ud305 \ucea0\ud398\uc778\uc740 `\"Munson's: More than Just Pickles and Preserves`\"\ub77c\ub294 \uc2ac\ub85c\uac74\uc744 \uc0ac\uc6a9\ud558\uc5ec Munson's \uc81c\ud488\uc758 \uc774\uc810\uc744 \uac15\uc870\ud569\ub2c8\ub2e4.\"\n  },\n  @{\n    Old = \"\ub9c8\ucf00\ud305 \ucea0\ud398\uc778\uc740 \ub2e4\uc74c \ubaa8\ud1a0\ub97c \uc0ac\uc6a9\ud558\uc5ec Munson\uc758 \uace0\uac1d \uc639\ud638\uc5d0 \uc601\uac10\uc744 \uc90d\ub2c8\ub2e4. `\"Munson's: \ud53c\ud074\uacfc \ubcf4\uc874\uc758 \uc0ac\ub791\uc744 \uacf5\uc720\ud558\uc2ed\uc2dc\uc624`\".\"\n    New = \"\ub9c8\ucf00\ud305 \ucea0\ud398\uc778\uc740 `\"Munson's: Share the Love of Pickles and Preserves`\"\ub77c\ub294 \ubaa8\ud1a0\ub97c \uc0ac\uc6a9\ud558\uc5ec Munson's \uace0\uac1d\uc758 \uc9c0\uc9c0\ub3c4\uc5d0 \uc601\uac10\uc744 \uc90d\ub2c8\ub2e4.\"\n  },\n  @{\n    Old = \"\ub9c8\ucf00\ud305 \ucea0\ud398\uc778\uc740 \ub2e4\uc74c \ubb38\uad6c\ub97c \uc0ac\uc6a9\ud558\uc5ec Munson\uc758 \uc81c\ud488 \ud3c9\uac00\ud310\uc744 \uad6c\ub3d9\ud558\uace0 \uad6c\ub9e4\ud569\ub2c8\ub2e4. `\"Munson's: \ucc3e\uae30, \uc0ac\uc6a9\ud574 \ubcf4\uae30, \uc0ac\ub791\ud558\uc138\uc694`\".\"\n    New = \"\ub9c8\ucf00\ud305 \ucea0\ud398\uc778\uc740 `\"Munson's: Find Them, Try Them, Love Them`\"\ub77c\ub294 \ubb38\uad6c\ub97c \uc0ac\uc6a9\ud558\uc5ec Munson's \uc81c\ud488\uc744 \uccb4\ud5d8\ud558\uace0 \uad6c\ub9e4\ud558\ub3c4\ub85d \ud569\ub2c8\ub2e4.\"\n  }\n)\n\nforeach ($rep in $replacements) {\n  # Build a fresh whole-document range for each search so earlier edits\n  # (which shift character offsets) cannot affect later lookups.\n  $rng = $d.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Text = $rep.Old\n  $rng.Find.Forward = $true\n  $rng.Find.Wrap = 1  # wdFindContinue\n\n  $found = $rng.Find.Execute()\n  if (-not $found) {\n    throw \"Text not found: $($rep.Old)\"\n  }\n\n  # A successful Execute() collapses $rng to the matched text, so setting\n  # .Text here replaces only that paragraph's run, preserving its\n  # surrounding run formatting (rPr) untouched.\n  $rng.Text = $rep.New\n}\n"}
